$d = $word.ActiveDocument
$d.Content.Find.Execute("away", $true, $false, $false, $false, $false, $true, 1, $false, "away,", 2)
